# Apply the "+1" counter updates produced by the latest data refresh.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 1523
$wsExpo.Range("F12").Value = 5282
$wsExpo.Range("F18").Value = 62
$wsExpo.Range("F24").Value = 3778

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 88

# Sheet "全部类型" (all types, aggregated)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 88
$wsAll.Range("F5").Value  = 1523
$wsAll.Range("F13").Value = 5282
$wsAll.Range("F19").Value = 62
$wsAll.Range("F25").Value = 3778
